$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last header cell (G1) to the new header cell (H1)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Add new "Save" header in column H
$ws.Range("H1").Value = "Save"

# Add the new numeric data values in column H for rows 2 and 3
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
